# Update the "generate report" timestamps that were refreshed when the
# handback status report was regenerated.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the
# d08eb7dc-77a2-4176-a7f0-6948c41606dc.md row
$wsOverview.Range("G4").Value = "2016-08-20 10:51:03"

# zh-cn sheet: Correspond Handoff / Handback datetimes for the same row
$wsZhCn.Range("H4").Value = "2016-08-20 10:50:57"
$wsZhCn.Range("K4").Value = "2016-08-20 10:51:27"

# de-de sheet: Correspond Handoff / Handback datetimes for the same row
$wsDeDe.Range("H4").Value = "2016-08-20 10:51:03"
$wsDeDe.Range("K4").Value = "2016-08-20 10:51:33"
